$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1193696.5
$ws.Range("I32").Value = 546075.9399999999
$ws.Range("K32").Value = 546075.9399999999
$ws.Range("M32").Value = -545788.9399999999

$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 28030.438
$ws.Range("I45").Value = 48310.223
$ws.Range("K45").Value = 48310.223
$ws.Range("M45").Value = -47933.223

$ws.Range("H61").Value = 1925.6333
$ws.Range("I61").Value = 1485.0454
$ws.Range("J61").Value = 3137.25
$ws.Range("K61").Value = 1485.0454
$ws.Range("L61").Value = 3137.25
$ws.Range("M61").Value = -1273.0454
$ws.Range("N61").Value = -3561.25

$ws.Range("H74").Value = 2125.182
$ws.Range("I74").Value = 1930.7059
$ws.Range("K74").Value = 1930.7059
$ws.Range("M74").Value = -1056.7059

$ws.Range("H77").Value = 2125.182
$ws.Range("I77").Value = 1930.7059
$ws.Range("K77").Value = 9653.529500000001
$ws.Range("M77").Value = -5285.529500000001

$ws.Range("H122").Value = 2567.5
$ws.Range("I122").Value = 2590.2
$ws.Range("K122").Value = 7770.599999999999
$ws.Range("M122").Value = -5320.599999999999

$ws.Range("H132").Value = 1980.8667
$ws.Range("I132").Value = 1148.7778
$ws.Range("K132").Value = 3446.3334
$ws.Range("M132").Value = -916.3334000000004

$ws.Range("H136").Value = 1925.6333
$ws.Range("I136").Value = 1485.0454
$ws.Range("J136").Value = 3137.25
$ws.Range("K136").Value = 4455.1362
$ws.Range("L136").Value = 9411.75
$ws.Range("M136").Value = -1905.1362
$ws.Range("N136").Value = -14511.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 11332.333
$ws.Range("I36").Value = 20000
$ws.Range("J36").Value = 9598.799999999999
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 9598.799999999999
$ws.Range("M36").Value = -19466
$ws.Range("N36").Value = -10666.8

$ws.Range("H82").Value = 27732.143
$ws.Range("I82").Value = 7532
$ws.Range("K82").Value = 7532
$ws.Range("M82").Value = -7149

$ws.Range("H85").Value = 27732.143
$ws.Range("I85").Value = 7532
$ws.Range("K85").Value = 7532
$ws.Range("M85").Value = -6206

$ws.Range("H86").Value = 5083.3335
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 5625
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 5625
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -7871

$ws.Range("H89").Value = 5083.3335
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 5625
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 28125
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -39357

$ws.Range("H99").Value = 3266.1667
$ws.Range("I99").Value = 2419.4
$ws.Range("K99").Value = 2419.4
$ws.Range("M99").Value = -921.4000000000001

$ws.Range("H139").Value = 67117.86
$ws.Range("J139").Value = 71640.5
$ws.Range("L139").Value = 71640.5
$ws.Range("N139").Value = -81920.5

$ws.Range("H141").Value = 78000
$ws.Range("J141").Value = 78000
$ws.Range("L141").Value = 78000
$ws.Range("N141").Value = -88360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 628000
$ws.Range("I4").Value = 1000000
$ws.Range("J4").Value = 535000
$ws.Range("K4").Value = 1000000
$ws.Range("L4").Value = 535000
$ws.Range("M4").Value = -999888
$ws.Range("N4").Value = -535224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 12508941
$ws.Range("I68").Value = 905.5
$ws.Range("K68").Value = 2716.5
$ws.Range("M68").Value = -1905.5

$ws.Range("H71").Value = 12508941
$ws.Range("I71").Value = 905.5
$ws.Range("K71").Value = 8149.5
$ws.Range("M71").Value = -4093.5

$ws.Range("H97").Value = 714773.1
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H141").Value = 21009.428
$ws.Range("I141").Value = 13413.2
$ws.Range("K141").Value = 40239.60000000001
$ws.Range("M141").Value = -35059.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16772.334
$ws.Range("I43").Value = 2658.5
$ws.Range("K43").Value = 2658.5
$ws.Range("M43").Value = -2507.5

$ws.Range("H132").Value = 2082.7837
$ws.Range("I132").Value = 1957.7333
$ws.Range("K132").Value = 5873.199900000001
$ws.Range("M132").Value = -3343.199900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1877.6666
$ws.Range("I7").Value = 1733.7
$ws.Range("J7").Value = 2597.5
$ws.Range("K7").Value = 1733.7
$ws.Range("L7").Value = 2597.5
$ws.Range("M7").Value = -1621.7
$ws.Range("N7").Value = -2821.5

$ws.Range("H16").Value = 2195.5715
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 3989.6667
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 3989.6667
$ws.Range("M16").Value = -680
$ws.Range("N16").Value = -4329.6667

$ws.Range("H40").Value = 24000.75
$ws.Range("I40").Value = 70004
$ws.Range("J40").Value = 8666.333000000001
$ws.Range("K40").Value = 70004
$ws.Range("L40").Value = 8666.333000000001
$ws.Range("M40").Value = -69868
$ws.Range("N40").Value = -8938.333000000001

$ws.Range("H46").Value = 4775
$ws.Range("J46").Value = 6250
$ws.Range("L46").Value = 6250
$ws.Range("N46").Value = -6626

$ws.Range("H122").Value = 9216
$ws.Range("J122").Value = 10982.134
$ws.Range("L122").Value = 32946.402
$ws.Range("N122").Value = -37846.402

$ws.Range("H126").Value = 1877.6666
$ws.Range("I126").Value = 1733.7
$ws.Range("J126").Value = 2597.5
$ws.Range("K126").Value = 5201.1
$ws.Range("L126").Value = 7792.5
$ws.Range("M126").Value = -2731.1
$ws.Range("N126").Value = -12732.5

$ws.Range("H132").Value = 13749.667
$ws.Range("I132").Value = 12999.667
$ws.Range("J132").Value = 14499.667
$ws.Range("K132").Value = 38999.001
$ws.Range("L132").Value = 43499.001
$ws.Range("M132").Value = -36469.001
$ws.Range("N132").Value = -48559.001

$ws.Range("H136").Value = 4938.625
$ws.Range("I136").Value = 2168.1667
$ws.Range("J136").Value = 13250
$ws.Range("K136").Value = 6504.500100000001
$ws.Range("L136").Value = 39750
$ws.Range("M136").Value = -3954.500100000001
$ws.Range("N136").Value = -44850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2333.3333
$ws.Range("I5").Value = 2333.3333
$ws.Range("K5").Value = 2333.3333
$ws.Range("M5").Value = -2221.3333

$ws.Range("H32").Value = 17304.4
$ws.Range("I32").Value = 12880.5
$ws.Range("K32").Value = 12880.5
$ws.Range("M32").Value = -12563.5

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H54").Value = 40880.625
$ws.Range("I54").Value = 22070
$ws.Range("J54").Value = 43567.855
$ws.Range("K54").Value = 22070
$ws.Range("L54").Value = 43567.855
$ws.Range("M54").Value = -21550
$ws.Range("N54").Value = -44607.855

$ws.Range("H62").Value = 6333
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 6333
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 3391.4583
$ws.Range("I81").Value = 2655.389
$ws.Range("J81").Value = 5599.6665
$ws.Range("K81").Value = 5310.778
$ws.Range("L81").Value = 11199.333
$ws.Range("M81").Value = -4249.778
$ws.Range("N81").Value = -13321.333

$ws.Range("H84").Value = 3391.4583
$ws.Range("I84").Value = 2655.389
$ws.Range("J84").Value = 5599.6665
$ws.Range("K84").Value = 26553.89
$ws.Range("L84").Value = 55996.665
$ws.Range("M84").Value = -21249.89
$ws.Range("N84").Value = -66604.66500000001

$ws.Range("H122").Value = 125002000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

$ws.Range("H132").Value = 3420.5652
$ws.Range("I132").Value = 3332.1667
$ws.Range("J132").Value = 3738.8
$ws.Range("K132").Value = 9996.500100000001
$ws.Range("L132").Value = 11216.4
$ws.Range("M132").Value = -7466.500100000001
$ws.Range("N132").Value = -16276.4

$ws.Range("H136").Value = 5923.7
$ws.Range("I136").Value = 4467.125
$ws.Range("J136").Value = 11750
$ws.Range("K136").Value = 13401.375
$ws.Range("L136").Value = 35250
$ws.Range("M136").Value = -10851.375
$ws.Range("N136").Value = -40350
